$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.324.61"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "1.987.36"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "'245.25"
$ws.Range("E5").Value = "  -2.76%  "

$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").Value = "'62.71"
$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("D10").Value = "'56.44"
$ws.Range("E10").Value = "  -4.06%  "

$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  +7.73%  "

$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").Value = "'0.870"
$ws.Range("E13").Value = "  -3.16%  "

$ws.Range("D14").Value = "'22.59"
$ws.Range("E14").Value = "  +11.64%  "

$ws.Range("D15").Value = "'14.12"
$ws.Range("E15").Value = "  -5.68%  "

$ws.Range("D16").Value = "2.279.36"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").Value = "'5.48"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "1.996.13"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D19").Value = "36.222.05"
$ws.Range("E19").Value = "  -2.16%  "

$ws.Range("D20").Value = "'71.43"
$ws.Range("E20").Value = "  -1.92%  "

$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").Value = "'5.30"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'237.55"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  -10.03%  "

$ws.Range("D26").Value = "'2.31"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("D27").Value = "'9.84"
$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  +23.04%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'160.12"
$ws.Range("E29").Value = "  -3.13%  "

$ws.Range("D30").Value = "'20.00"
$ws.Range("E30").Value = "  +1.31%  "

$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("E32").Value = "  -3.82%  "

$ws.Range("E33").Value = "  -4.59%  "

$ws.Range("E34").Value = "  +2.50%  "

$ws.Range("E35").Value = "  -5.63%  "

$ws.Range("D36").Value = "'6.38"
$ws.Range("E36").Value = "  +6.78%  "

$ws.Range("D37").Value = "'2.29"
$ws.Range("E37").Value = "  -5.77%  "

$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("D40").Value = "'3.15"
$ws.Range("E40").Value = "  +15.44%  "

$ws.Range("D41").Value = "'0.0993"
$ws.Range("E41").Value = "  -4.77%  "

$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.87"
$ws.Range("E43").Value = "  -2.06%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0215"
$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("E45").Value = "  -2.63%  "

$ws.Range("D46").Value = "'93.60"
$ws.Range("E46").Value = "  -1.22%  "

$ws.Range("D47").Value = "'16.33"
$ws.Range("E47").Value = "  -2.62%  "

$ws.Range("D48").Value = "'7.56"
$ws.Range("E48").Value = "  -5.73%  "

$ws.Range("D49").Value = "1.354.30"
$ws.Range("E49").Value = "  -4.91%  "

$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("D51").Value = "2.173.33"
$ws.Range("E51").Value = "  -1.79%  "

